$d = $word.ActiveDocument

# The run "<id>", the run "p056v_1", and the run "</id>" (three separately
# formatted runs) need to become a single run containing the full text
# "<id>p056v_1</id>", using the formatting of the first ("<id>") run.
# A find/replace across the run boundary merges the matched text into one
# run that inherits the formatting of the run where the match starts.
$d.Content.Find.Execute("<id>p056v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p056v_1</id>", 2)
